$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name text (shared by cell B1 on both sheets) to the new name.
$newProductName = "4246-RBI-EI-DB-SAR-REC-RNI-INT-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-PER-1st"
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Change the short name value from the numeric 4246 to the text "424s" to remove
# the test case inter-dependency.
$wsInput.Range("B2").Value = "424s"

# Update the saved selection on the input sheet so only B2 (not B2:B3) is selected,
# then restore the originally active sheet (ProductLoanOutput) so that doesn't change.
[void]$wsInput.Range("B2").Select()
$wsOutput.Activate()
